# Updates the "INCO" sheet of mapa_interactivo_INCO.xlsx:
#   - inserts 4 new incident rows before the old last row (46-49)
#   - the old last row (previously 46) shifts down to row 50, unchanged
#   - appends 4 more new incident rows at the end (51-54), with blank
#     coordinate/operation/zone columns (not geolocated yet)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the previous last data row (46) down to 50, opening up 4 fresh rows
$ws.Rows("46:49").Insert()

# --- Row 46: Caso -539 ---
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = '-539'
$ws.Range("A46").Style = "Normal"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = '7/31/2025'
$ws.Range("B46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '808615951'
$ws.Range("E46").Style = "Normal"
$ws.Range("C46").Value = 'Tejedor 1097'
$ws.Range("F46").Value = 'INCO'
$ws.Range("G46").Value = 'Pendiente'
$ws.Range("H46").Value = 'Picada'
$ws.Range("J46").Value = 'Cambio'
$ws.Range("K46").Value = 'Sin equipos'
$ws.Range("L46").Value = 'Terminal'
$ws.Range("I46").Value = 1
$ws.Range("M46").Value = -58.440748
$ws.Range("N46").Value = -34.63245
$ws.Range("O46").Value = 'Boedo'
$ws.Range("P46").Value = 'Capital Sur'

# --- Row 47: Caso -540 ---
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = '-540'
$ws.Range("A47").Style = "Normal"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = '7/31/2025'
$ws.Range("B47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '808615948'
$ws.Range("E47").Style = "Normal"
$ws.Range("C47").Value = 'Tejedor 1071'
$ws.Range("F47").Value = 'INCO'
$ws.Range("G47").Value = 'Pendiente'
$ws.Range("H47").Value = 'Picada'
$ws.Range("J47").Value = 'Cambio'
$ws.Range("K47").Value = 'Sin equipos'
$ws.Range("L47").Value = 'Pasante'
$ws.Range("I47").Value = 1
$ws.Range("M47").Value = -58.44037
$ws.Range("N47").Value = -34.632249
$ws.Range("O47").Value = 'Boedo'
$ws.Range("P47").Value = 'Capital Sur'

# --- Row 48: Caso -542 ---
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = '-542'
$ws.Range("A48").Style = "Normal"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = '8/1/2025'
$ws.Range("B48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '808663881'
$ws.Range("E48").Style = "Normal"
$ws.Range("C48").Value = 'Cramer 2141'
$ws.Range("F48").Value = 'INCO'
$ws.Range("G48").Value = 'Pendiente'
$ws.Range("H48").Value = 'Cambiar columna 114 base corroida '
$ws.Range("J48").Value = 'Cambio'
$ws.Range("K48").Value = 'Sin equipos'
$ws.Range("L48").Value = 'Pasante'
$ws.Range("I48").Value = 1
$ws.Range("M48").Value = -58.461582
$ws.Range("N48").Value = -34.564296
$ws.Range("O48").Value = 'Saavedra'
$ws.Range("P48").Value = 'Capital Norte'

# --- Row 49: Caso -544 ---
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = '-544'
$ws.Range("A49").Style = "Normal"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = '8/2/2025'
$ws.Range("B49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '15'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '808669129'
$ws.Range("E49").Style = "Normal"
$ws.Range("C49").Value = 'Vera 453'
$ws.Range("F49").Value = 'INCO'
$ws.Range("G49").Value = 'Pendiente'
$ws.Range("H49").Value = 'Columna corroída en base'
$ws.Range("J49").Value = 'Cambio'
$ws.Range("K49").Value = 'Sin equipos'
$ws.Range("L49").Value = 'Pasante'
$ws.Range("I49").Value = 1
$ws.Range("M49").Value = -58.437239
$ws.Range("N49").Value = -34.599438
$ws.Range("O49").Value = 'Palermo'
$ws.Range("P49").Value = 'Capital Sur'

# --- Row 51: Caso -552 ---
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = '-552'
$ws.Range("A51").Style = "Normal"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = '8/14/2025'
$ws.Range("B51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '808973183'
$ws.Range("E51").Style = "Normal"
$ws.Range("C51").Value = 'Catulo Castillo 2890'
$ws.Range("F51").Value = 'INCO'
$ws.Range("G51").Value = 'Pendiente'
$ws.Range("H51").Value = 'Picada'
$ws.Range("J51").Value = 'Cambio'
$ws.Range("K51").Value = 'Sin equipos'
$ws.Range("L51").Value = 'Pasante'
$ws.Range("I51").Value = 1
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = ""
$ws.Range("O51").Value = ""
$ws.Range("P51").Value = ""

# --- Row 52: Caso -553 ---
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = '-553'
$ws.Range("A52").Style = "Normal"
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = '8/14/2025'
$ws.Range("B52").Style = "Normal"
$ws.Range("D52").NumberFormat = "@"
$ws.Range("D52").Value = '12'
$ws.Range("D52").Style = "Normal"
$ws.Range("E52").NumberFormat = "@"
$ws.Range("E52").Value = '808973192'
$ws.Range("E52").Style = "Normal"
$ws.Range("C52").Value = 'Holmberg 4002'
$ws.Range("F52").Value = 'INCO'
$ws.Range("G52").Value = 'Pendiente'
$ws.Range("H52").Value = 'Edificio en construcción solicitan correr columna 114 por entrada de garaje 5mts aprox Hablar con Sr Galvan encargado de la Obra'
$ws.Range("J52").Value = 'Cambio'
$ws.Range("K52").Value = 'Sin equipos'
$ws.Range("L52").Value = 'Pasante'
$ws.Range("I52").Value = 1
$ws.Range("M52").Value = ""
$ws.Range("N52").Value = ""
$ws.Range("O52").Value = ""
$ws.Range("P52").Value = ""

# --- Row 53: Caso -554 ---
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = '-554'
$ws.Range("A53").Style = "Normal"
$ws.Range("B53").NumberFormat = "@"
$ws.Range("B53").Value = '8/14/2025'
$ws.Range("B53").Style = "Normal"
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = '1'
$ws.Range("D53").Style = "Normal"
$ws.Range("E53").NumberFormat = "@"
$ws.Range("E53").Value = '808973197'
$ws.Range("E53").Style = "Normal"
$ws.Range("C53").Value = 'Lima Oeste 1697'
$ws.Range("F53").Value = 'INCO'
$ws.Range("G53").Value = 'Pendiente'
$ws.Range("H53").Value = 'Columna inclinada'
$ws.Range("J53").Value = 'Aplomo'
$ws.Range("K53").Value = 'Sin equipos'
$ws.Range("L53").Value = 'Terminal'
$ws.Range("I53").Value = 1
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = ""
$ws.Range("O53").Value = ""
$ws.Range("P53").Value = ""

# --- Row 54: Caso -555 ---
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value = '-555'
$ws.Range("A54").Style = "Normal"
$ws.Range("B54").NumberFormat = "@"
$ws.Range("B54").Value = '8/14/2025'
$ws.Range("B54").Style = "Normal"
$ws.Range("D54").NumberFormat = "@"
$ws.Range("D54").Value = '1'
$ws.Range("D54").Style = "Normal"
$ws.Range("E54").NumberFormat = "@"
$ws.Range("E54").Value = '808973201'
$ws.Range("E54").Style = "Normal"
$ws.Range("C54").Value = 'Lima Oeste 1649'
$ws.Range("F54").Value = 'INCO'
$ws.Range("G54").Value = 'Pendiente'
$ws.Range("H54").Value = 'Picada'
$ws.Range("J54").Value = 'Cambio'
$ws.Range("K54").Value = 'Sin equipos'
$ws.Range("L54").Value = 'Pasante'
$ws.Range("I54").Value = 1
$ws.Range("M54").Value = ""
$ws.Range("N54").Value = ""
$ws.Range("O54").Value = ""
$ws.Range("P54").Value = ""

Write-Output "mapa_interactivo_INCO: added rows 46-54 (dimension now A1:P54)"
